$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated FilesTab query (B4): removed the "File Type" and "Breed" output columns
$newFileQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN['Black and Tan Coonhound']`nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFileQuery

# Row 4 height shrinks because the query text is shorter now
$ws.Rows.Item(4).RowHeight = 217.5

# Update the active selection to C4 (matches the recorded view state after the edit)
$ws.Range("C4").Select()
